$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 203
$ws.Range("C2").Value = 70.23999999999999

$ws.Range("B3").Value = 86
$ws.Range("C3").Value = 29.76
